# Insert a new section row for course "BIOL222" above the CHEM110 rows.
# This pushes the existing row 9 ("CHEM110", ...) and everything below it
# down by one row (old row N -> new row N+1 for N >= 9), and fills the
# newly-inserted row 9 with the BIOL222 section data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at row 9 (shifts CHEM110... rows down to 10...45).
$ws.Rows.Item(9).Insert()

# Populate the new row with the BIOL222 lab section.
$ws.Range("A9").Value = "BIOL222"
$ws.Range("B9").Value = 1
$ws.Range("C9").Value = 1
$ws.Range("D9").Value = "Wednesday"
$ws.Range("E9").Value = 0.39930555555555558
$ws.Range("F9").Value = 0.125
$ws.Range("G9").Value = 1

# Match the author's final selection in the saved workbook.
$ws.Range("A10").Select()
